# Updates cryptos list values (Price / Volume(1h) columns) per latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.412.04'
$ws.Range("E2").Value = '  +4.41%  '

$ws.Range("D3").Value = '2.045.05'
$ws.Range("E3").Value = '  +2.99%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = '''252.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.92%  '

$ws.Range("D6").Value = '''0.651'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.18%  '

$ws.Range("D7").Value = '''66.23'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +11.02%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '''0.400'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +9.69%  '

$ws.Range("D10").Value = '''59.21'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.03%  '

$ws.Range("D11").Value = '''0.0799'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.15%  '

$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '''0.915'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.98%  '

$ws.Range("D14").Value = '''23.57'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +23.86%  '

$ws.Range("D15").Value = '''14.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.20%  '

$ws.Range("D16").Value = '2.343.90'
$ws.Range("E16").Value = '  +2.94%  '

$ws.Range("D17").Value = '''5.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +7.67%  '

$ws.Range("D18").Value = '2.040.59'
$ws.Range("E18").Value = '  +2.68%  '

$ws.Range("D19").Value = '37.282.21'
$ws.Range("E19").Value = '  +4.27%  '

$ws.Range("D20").Value = '''73.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.30%  '

$ws.Range("D21").Value = '0.0₃0896'
$ws.Range("E21").Value = '  +5.84%  '

$ws.Range("D22").Value = '''5.54'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.68%  '

$ws.Range("D23").Value = '''239.26'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.88%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  -0.75%  '

$ws.Range("D26").Value = '''2.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.38%  '

$ws.Range("D27").Value = '''10.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.03%  '

$ws.Range("D28").Value = '''161.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.75%  '

$ws.Range("D29").Value = '''20.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.72%  '

$ws.Range("D30").Value = '''0.129'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +32.58%  '

$ws.Range("D31").Value = '''0.123'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.12%  '

$ws.Range("D32").Value = '''5.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.78%  '

$ws.Range("E33").Value = '  +5.52%  '

$ws.Range("D34").Value = '''0.0631'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.03%  '

$ws.Range("D35").Value = '''4.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.96%  '

$ws.Range("D36").Value = '''6.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +12.70%  '

$ws.Range("E37").Value = '  -2.78%  '

$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("E39").Value = '  +3.74%  '

$ws.Range("D40").Value = '''2.99'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +29.34%  '

$ws.Range("D41").Value = '''1.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.81%  '

$ws.Range("D42").Value = '''0.101'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.29%  '

$ws.Range("E43").Value = '  +5.75%  '

$ws.Range("E44").Value = '  +6.27%  '

$ws.Range("D45").Value = '''17.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.60%  '

$ws.Range("E46").Value = '  +2.79%  '

$ws.Range("D47").Value = '''95.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.48%  '

$ws.Range("E48").Value = '  +1.52%  '

$ws.Range("D49").Value = '1.395.90'
$ws.Range("E49").Value = '  +2.50%  '

$ws.Range("E50").Value = '  +0.74%  '

$ws.Range("D51").Value = '''47.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.26%  '
